# Update the SQL query text for the "Website" field in the ProgramsTab StatQuery (cell B2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "SELECT DISTINCT `n    prg.program_name AS ""Program"",`n  CASE`n    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym`n        ELSE prg.program_link`n    END  AS ""Website"",`n    prg.focus_area AS ""Focus Area"",`n    prg.cancer_type AS ""Cancer Type"",`n CASE `n        WHEN prg.data_link IS NOT NULL THEN prg.website       `n        ELSE prg.data_link`n    END AS ""Data Location Details""`nFROM `n    df_program prg`nWHERE `n     prg.cancer_type LIKE '%Lymphoma%'`nORDER BY `n    lower(prg.program_name) ASC`nLIMIT 100;"

$b2 = $ws.Range("B2")

# Re-touch the font formatting so the cell is re-written with its own (new) style entry,
# matching the re-formatting that happened when this query text was edited/pasted in Excel.
$b2.Font.ThemeColor = 1
$b2.Font.Size = 12
$b2.WrapText = $true

$b2.Value = $newText

# Reset the sheet view: scroll back to the top-left (removes the old topLeftCell="A5")
# and move the active selection to B2 (previously it was left on B11).
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B2").Select() | Out-Null
